$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 02:08"

# Re-rank country name labels that swapped position in the sorted table
$ws.Range("A107").Value = "Luxemburgo"
$ws.Range("A108").Value = "Mauritania"
$ws.Range("A118").Value = "Cabo Verde"
$ws.Range("A119").Value = "Nicaragua"
$ws.Range("A120").Value = "Congo"
$ws.Range("A123").Value = "Surinam"
$ws.Range("A124").Value = "Ruanda"
$ws.Range("A129").Value = "Angola"
$ws.Range("A130").Value = "Siria"

# Refresh updated per-country statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes)
$ws.Range("B4").Value = 6825967
$ws.Range("C4").Value = 37820
$ws.Range("D4").Value = 4103090
$ws.Range("E4").Value = 2521599
$ws.Range("G4").Value = 1081
$ws.Range("H4").Value = 201278

$ws.Range("B6").Value = 4421686
$ws.Range("C6").Value = 37387
$ws.Range("E6").Value = 567200
$ws.Range("G6").Value = 967
$ws.Range("H6").Value = 134174

$ws.Range("B13").Value = 589012
$ws.Range("C13").Value = 11674
$ws.Range("E13").Value = 128633
$ws.Range("G13").Value = 264
$ws.Range("H13").Value = 12116

$ws.Range("B36").Value = 103466
$ws.Range("C36").Value = 634
$ws.Range("D36").Value = 76787
$ws.Range("E36").Value = 24481
$ws.Range("G36").Value = 11
$ws.Range("H36").Value = 2198

$ws.Range("D57").Value = 56955
$ws.Range("E57").Value = 532

$ws.Range("B66").Value = 41032
$ws.Range("C66").Value = 2136
$ws.Range("D66").Value = 22931
$ws.Range("E66").Value = 17619
$ws.Range("G66").Value = 6
$ws.Range("H66").Value = 482

$ws.Range("B103").Value = 8678
$ws.Range("C103").Value = 24
$ws.Range("D103").Value = 7827
$ws.Range("E103").Value = 798

$ws.Range("B104").Value = 8541
$ws.Range("C104").Value = 11
$ws.Range("D104").Value = 6258
$ws.Range("E104").Value = 2063

$ws.Range("B105").Value = 8100
$ws.Range("C105").Value = 477
$ws.Range("D105").Value = 2309
$ws.Range("E105").Value = 5662
$ws.Range("G105").Value = 6
$ws.Range("H105").Value = 129

$ws.Range("B106").Value = 7598
$ws.Range("C106").Value = 22
$ws.Range("D106").Value = 5823
$ws.Range("E106").Value = 1551

$ws.Range("B107").Value = 7394
$ws.Range("C107").Value = 110
$ws.Range("D107").Value = 6593
$ws.Range("E107").Value = 677
$ws.Range("H107").Value = 124

$ws.Range("B108").Value = 7332
$ws.Range("C108").Value = 13
$ws.Range("D108").Value = 6844
$ws.Range("E108").Value = 327
$ws.Range("H108").Value = 161

$ws.Range("B109").Value = 7291
$ws.Range("C109").Value = 230
$ws.Range("D109").Value = 4764
$ws.Range("E109").Value = 2401
$ws.Range("G109").Value = 3
$ws.Range("H109").Value = 126

$ws.Range("B118").Value = 4978
$ws.Range("C118").Value = 74
$ws.Range("D118").Value = 4430
$ws.Range("E118").Value = 501
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 47

$ws.Range("B119").Value = 4961
$ws.Range("D119").Value = 2913
$ws.Range("E119").Value = 1901
$ws.Range("H119").Value = 147

$ws.Range("B120").Value = 4934
$ws.Range("D120").Value = 3887
$ws.Range("E120").Value = 959
$ws.Range("H120").Value = 88

$ws.Range("B122").Value = 4782
$ws.Range("C122").Value = 7
$ws.Range("D122").Value = 1830
$ws.Range("E122").Value = 2890

$ws.Range("B123").Value = 4645
$ws.Range("C123").Value = 20
$ws.Range("D123").Value = 4089
$ws.Range("E123").Value = 461
$ws.Range("H123").Value = 95

$ws.Range("B124").Value = 4634
$ws.Range("C124").Value = 10
$ws.Range("D124").Value = 2789
$ws.Range("E124").Value = 1823
$ws.Range("H124").Value = 22

$ws.Range("B129").Value = 3675
$ws.Range("C129").Value = 106
$ws.Range("D129").Value = 1401
$ws.Range("E129").Value = 2131
$ws.Range("G129").Value = 4
$ws.Range("H129").Value = 143

$ws.Range("B130").Value = 3654
$ws.Range("C130").Value = 40
$ws.Range("D130").Value = 889
$ws.Range("E130").Value = 2602
$ws.Range("G130").Value = 3
$ws.Range("H130").Value = 163

$ws.Range("B154").Value = 1856
$ws.Range("C154").Value = 29
$ws.Range("D154").Value = 1559
$ws.Range("E154").Value = 252

$ws.Range("B156").Value = 1748
$ws.Range("C156").Value = 15
$ws.Range("D156").Value = 1162
$ws.Range("E156").Value = 530

$ws.Range("B158").Value = 1548
$ws.Range("C158").Value = 8
$ws.Range("E158").Value = 244

$ws.Range("B190").Value = 185
$ws.Range("C190").Value = 1
$ws.Range("D190").Value = 171
$ws.Range("E190").Value = 7

$ws.Range("B192").Value = 178
$ws.Range("C192").Value = 1
$ws.Range("E192").Value = 5
